# Commit: "Update ilr by ssa data"
#
# The Individualised Learner Record (ILR) source referenced by the FE
# achievements / participation / starts metrics moves on from the
# AY21/22 release to the AY22/23 release: the "LatestPeriod" label and
# the "sourceText" hyperlink both need updating for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

$oldPeriod = "AY21/22 data"
$newPeriod = "AY22/23 data"
$oldLink = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/further-education-and-skills/2021-22'>Individualised Learner Record</a>"
$newLink = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/further-education-and-skills/2022-23'>Individualised Learner Record</a>"

# Find the extent of the used range on the dataText sheet.
$lastRow = $ws.Cells.SpecialCells(11).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $periodCell = $ws.Cells.Item($r, 2)
    $linkCell = $ws.Cells.Item($r, 4)

    if ($periodCell.Value2 -eq $oldPeriod -and $linkCell.Value2 -eq $oldLink) {
        $periodCell.Value = $newPeriod
        $linkCell.Value = $newLink
    }
}

# The author's selection moved from F12 to D12 when saving.
[void]$ws.Range("D12").Select()
